# Loan RBI, Variable Instalments
#
# On the "Repayment Schedule" sheet a new column is inserted before column N
# (shifting the existing "Late"/heading/"Outstanding" columns one place to
# the right), and the "Repayment Schedule" tab becomes the active sheet with
# cell K14 selected (previously "Transactions" was the active tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N, pushing the old N/O/P columns
# to O/P/Q respectively.
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet and select K14, matching the
# new sheetView/tabSelected/activeTab state captured in the workbook.
$ws.Activate()
$ws.Range("K14").Select()
